$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update industry name in A16 (shared string "Banks - Regional " -> "Drug Manufacturers - Specialty Generic ")
$ws.Range("A16").Value = "Drug Manufacturers - Specialty Generic "

# Update frequency counts
$ws.Range("B5").Value = 2
$ws.Range("B7").Value = 3
$ws.Range("B8").Value = 3
